$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1222.1666
$ws.Range("J19").Value = 1252.1111
$ws.Range("L19").Value = 1252.1111
$ws.Range("N19").Value = -1602.1111

$ws.Range("H32").Value = 1390.3636
$ws.Range("J32").Value = 1307.25
$ws.Range("L32").Value = 1307.25
$ws.Range("N32").Value = -1959.25

$ws.Range("H87").Value = 71459.89
$ws.Range("J87").Value = 73604.586
$ws.Range("L87").Value = 73604.586
$ws.Range("N87").Value = -76100.586

$ws.Range("H90").Value = 71459.89
$ws.Range("J90").Value = 73604.586
$ws.Range("L90").Value = 220813.758
$ws.Range("N90").Value = -233293.758

$ws.Range("H112").Value = 4599.92
$ws.Range("J112").Value = 4749.9165
$ws.Range("L112").Value = 14249.7495
$ws.Range("N112").Value = -16465.7495

$ws.Range("H127").Value = 2389
$ws.Range("I127").Value = 961.46155
$ws.Range("J127").Value = 3626.2
$ws.Range("K127").Value = 2884.38465
$ws.Range("L127").Value = 10878.6
$ws.Range("M127").Value = 2075.61535
$ws.Range("N127").Value = -20798.6

$ws.Range("H132").Value = 121177.39
$ws.Range("I132").Value = 319742.7
$ws.Range("J132").Value = 14257.615
$ws.Range("K132").Value = 959228.1000000001
$ws.Range("L132").Value = 42772.845
$ws.Range("M132").Value = -956698.1000000001
$ws.Range("N132").Value = -47832.845

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 14078448
$ws.Range("I2").Value = 1343754.5
$ws.Range("J2").Value = 41670284
$ws.Range("K2").Value = 1343754.5
$ws.Range("L2").Value = 41670284
$ws.Range("M2").Value = -1343641.5
$ws.Range("N2").Value = -41670510

$ws.Range("H61").Value = 12122.272
$ws.Range("I61").Value = 12790.889
$ws.Range("J61").Value = 11319.934
$ws.Range("K61").Value = 12790.889
$ws.Range("L61").Value = 11319.934
$ws.Range("M61").Value = -12578.889
$ws.Range("N61").Value = -11743.934

$ws.Range("H74").Value = 5437712
$ws.Range("I74").Value = 10000962
$ws.Range("K74").Value = 10000962
$ws.Range("M74").Value = -10000088

$ws.Range("H77").Value = 5437712
$ws.Range("I77").Value = 10000962
$ws.Range("K77").Value = 50004810
$ws.Range("M77").Value = -50000442

$ws.Range("H116").Value = 14078448
$ws.Range("I116").Value = 1343754.5
$ws.Range("J116").Value = 41670284
$ws.Range("K116").Value = 1343754.5
$ws.Range("L116").Value = 41670284
$ws.Range("M116").Value = -1341460.5
$ws.Range("N116").Value = -41674872

$ws.Range("H132").Value = 4141.7656
$ws.Range("I132").Value = 3193
$ws.Range("K132").Value = 9579
$ws.Range("M132").Value = -7049

$ws.Range("H136").Value = 12122.272
$ws.Range("I136").Value = 12790.889
$ws.Range("J136").Value = 11319.934
$ws.Range("K136").Value = 38372.667
$ws.Range("L136").Value = 33959.802
$ws.Range("M136").Value = -35822.667
$ws.Range("N136").Value = -39059.802

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 14078448
$ws.Range("I3").Value = 1343754.5
$ws.Range("J3").Value = 41670284
$ws.Range("K3").Value = 1343754.5
$ws.Range("L3").Value = 41670284
$ws.Range("M3").Value = -1343640.5
$ws.Range("N3").Value = -41670512

$ws.Range("H20").Value = 1862.8
$ws.Range("I20").Value = 1204.0714
$ws.Range("K20").Value = 1204.0714
$ws.Range("M20").Value = -957.0714

$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").ClearContents()

$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").ClearContents()

$ws.Range("H99").Value = 26867870
$ws.Range("I99").Value = 68057610
$ws.Range("J99").Value = 4999.7393
$ws.Range("K99").Value = 68057610
$ws.Range("L99").Value = 4999.7393
$ws.Range("M99").Value = -68056112
$ws.Range("N99").Value = -7995.7393

$ws.Range("H134").Value = 3055.4143
$ws.Range("I134").Value = 2600.0852
$ws.Range("J134").Value = 3985.8696
$ws.Range("K134").Value = 7800.2556
$ws.Range("L134").Value = 11957.6088
$ws.Range("M134").Value = -5265.2556
$ws.Range("N134").Value = -17027.6088

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("I16").Value = 1759.625
$ws.Range("J16").Value = 786.25
$ws.Range("K16").Value = 1759.625
$ws.Range("L16").Value = 786.25
$ws.Range("M16").Value = -1472.625
$ws.Range("N16").Value = -1360.25

$ws.Range("H31").Value = 30307958
$ws.Range("I31").Value = 83334760
$ws.Range("J31").Value = 6927.1904
$ws.Range("K31").Value = 83334760
$ws.Range("L31").Value = 6927.1904
$ws.Range("M31").Value = -83334465
$ws.Range("N31").Value = -7517.1904

$ws.Range("H34").Value = 30307958
$ws.Range("I34").Value = 83334760
$ws.Range("J34").Value = 6927.1904
$ws.Range("K34").Value = 83334760
$ws.Range("L34").Value = 6927.1904
$ws.Range("M34").Value = -83334558
$ws.Range("N34").Value = -7331.1904

$ws.Range("H58").Value = 2149.5715
$ws.Range("I58").Value = 1887.25
$ws.Range("K58").Value = 1887.25
$ws.Range("M58").Value = -1684.25

$ws.Range("H75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").ClearContents()

$ws.Range("H78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").ClearContents()

$ws.Range("I94").Value = 182.5
$ws.Range("J94").Value = 1712.375
$ws.Range("K94").Value = 182.5
$ws.Range("L94").Value = 1712.375
$ws.Range("M94").Value = 268.5
$ws.Range("N94").Value = -2614.375

$ws.Range("I113").Value = 1759.625
$ws.Range("J113").Value = 786.25
$ws.Range("K113").Value = 1759.625
$ws.Range("L113").Value = 786.25
$ws.Range("M113").Value = 410.375
$ws.Range("N113").Value = -5126.25

$ws.Range("H122").Value = 3019.9546
$ws.Range("I122").Value = 1115.3
$ws.Range("K122").Value = 3345.9
$ws.Range("M122").Value = -895.8999999999996

$ws.Range("H136").Value = 2149.5715
$ws.Range("I136").Value = 1887.25
$ws.Range("K136").Value = 5661.75
$ws.Range("M136").Value = -3111.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 0
$ws.Range("M17").ClearContents()
$ws.Range("N17").ClearContents()

$ws.Range("H95").Value = 15000
$ws.Range("J95").Value = 15000
$ws.Range("L95").Value = 45000
$ws.Range("N95").Value = -49118

$ws.Range("H98").Value = 930
$ws.Range("J98").Value = 999.5238000000001
$ws.Range("L98").Value = 2998.5714
$ws.Range("N98").Value = -5994.571400000001

$ws.Range("H107").Value = 994193.1
$ws.Range("I107").Value = 1810
$ws.Range("K107").Value = 5430
$ws.Range("M107").Value = -3510

$ws.Range("H121").Value = 164246.38
$ws.Range("I121").Value = 325242.25
$ws.Range("J121").Value = 3250.5
$ws.Range("K121").Value = 975726.75
$ws.Range("L121").Value = 9751.5
$ws.Range("M121").Value = -974416.75
$ws.Range("N121").Value = -12371.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2957.1707
$ws.Range("I126").Value = 2013.091
$ws.Range("K126").Value = 6039.272999999999
$ws.Range("M126").Value = -3569.272999999999

$ws.Range("H132").Value = 2831.838
$ws.Range("I132").Value = 2638.8572
$ws.Range("J132").Value = 3085.125
$ws.Range("K132").Value = 7916.571599999999
$ws.Range("L132").Value = 9255.375
$ws.Range("M132").Value = -5386.571599999999
$ws.Range("N132").Value = -14315.375

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 27087898
$ws.Range("I40").Value = 7696695
$ws.Range("K40").Value = 7696695
$ws.Range("M40").Value = -7696559

$ws.Range("H76").Value = 31974.334
$ws.Range("J76").Value = 31974.334
$ws.Range("L76").Value = 31974.334
$ws.Range("N76").Value = -32650.334

$ws.Range("H79").Value = 31974.334
$ws.Range("J79").Value = 31974.334
$ws.Range("L79").Value = 31974.334
$ws.Range("N79").Value = -34314.334

$ws.Range("H93").Value = 3234.4
$ws.Range("I93").Value = 3234.4
$ws.Range("K93").Value = 3234.4
$ws.Range("M93").Value = -1986.4

$ws.Range("H132").Value = 5457.494
$ws.Range("I132").Value = 4884.0234
$ws.Range("J132").Value = 6073.975
$ws.Range("K132").Value = 14652.0702
$ws.Range("L132").Value = 18221.925
$ws.Range("M132").Value = -12122.0702
$ws.Range("N132").Value = -23281.925

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2471.5112
$ws.Range("I122").Value = 2300.4092
$ws.Range("K122").Value = 6901.2276
$ws.Range("M122").Value = -4451.2276

$ws.Range("H136").Value = 9422.958000000001
$ws.Range("I136").Value = 2100
$ws.Range("J136").Value = 9969.448
$ws.Range("K136").Value = 6300
$ws.Range("L136").Value = 29908.344
$ws.Range("M136").Value = -3750
$ws.Range("N136").Value = -35008.344
